# Insert a new weekly price record as row 94, pushing the existing
# rows 94:145 down to 95:146 (dimension grows from A1:R145 to A1:R146).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 94 downwards by inserting a new row above the current row 94.
# This preserves formatting (e.g. the date style on column D) the same
# way Excel's own "Insert Row" command does.
$ws.Rows(94).Insert()

# Populate the newly inserted row 94 with the new record.
$ws.Cells.Item(94, 1).Value = 7
$ws.Cells.Item(94, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(94, 3).Value = "Ñuble"
$ws.Cells.Item(94, 4).Value = 45097
$ws.Cells.Item(94, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(94, 5).Value = 16
$ws.Cells.Item(94, 6).Value = 100112031
$ws.Cells.Item(94, 7).Value = "Poroto verde"
$ws.Cells.Item(94, 8).Value = "Magnum"
$ws.Cells.Item(94, 9).Value = "Primera"
$ws.Cells.Item(94, 10).Value = 30
$ws.Cells.Item(94, 11).Value = 25000
$ws.Cells.Item(94, 12).Value = 25000
$ws.Cells.Item(94, 13).Value = 25000
$ws.Cells.Item(94, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(94, 15).Value = "Perú"
$ws.Cells.Item(94, 16).Value = 1000
$ws.Cells.Item(94, 17).Value = 25
$ws.Cells.Item(94, 18).Value = "Hortaliza"
